$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 68.76957566666668
$ws.Range("H2").Value = 206.308727
$ws.Range("I2").Value = 0.08013720535972996
$ws.Range("J2").Value = 0.08013720535972994
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 10967.62206518989
$ws.Range("R2").Value = 98708.59858670905
$ws.Range("S2").Value = 0.02390721294954317
$ws.Range("T2").Value = 0.02390721294954316

$ws.Range("G3").Value = 68.76957566666668
$ws.Range("H3").Value = 206.308727
$ws.Range("I3").Value = 0.08013720535972996
$ws.Range("J3").Value = 0.08013720535972994
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 11866.77461636777
$ws.Range("R3").Value = 106800.97154731
$ws.Range("S3").Value = 0.02586718489126073
$ws.Range("T3").Value = 0.02586718489126072

$ws.Range("G4").Value = 68.76957566666668
$ws.Range("H4").Value = 206.308727
$ws.Range("I4").Value = 0.08013720535972996
$ws.Range("J4").Value = 0.08013720535972994
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 5115.610953513563
$ws.Range("R4").Value = 46040.49858162207
$ws.Range("S4").Value = 0.01115100426562217
$ws.Range("T4").Value = 0.01115100426562217

$ws.Range("G5").Value = 68.76957566666668
$ws.Range("H5").Value = 206.308727
$ws.Range("I5").Value = 0.08013720535972996
$ws.Range("J5").Value = 0.08013720535972994
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 4017.148240435319
$ws.Range("R5").Value = 36154.33416391787
$ws.Range("S5").Value = 0.008756576207962812
$ws.Range("T5").Value = 0.008756576207962812

$ws.Range("G6").Value = 68.76957566666668
$ws.Range("H6").Value = 206.308727
$ws.Range("I6").Value = 0.08013720535972996
$ws.Range("J6").Value = 0.08013720535972994
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 4796.417678675679
$ws.Range("R6").Value = 43167.7591080811
$ws.Range("S6").Value = 0.01045522704534108
$ws.Range("T6").Value = 0.01045522704534108

$ws.Range("I7").Value = 0.2583188840501616
$ws.Range("J7").Value = 0.2583188840501616
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 35353.66475341897
$ws.Range("R7").Value = 318182.9827807707
$ws.Range("S7").Value = 0.07706388739354429
$ws.Range("T7").Value = 0.07706388739354429

$ws.Range("I8").Value = 0.2583188840501616
$ws.Range("J8").Value = 0.2583188840501616
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.0833817738544132
$ws.Range("T8").Value = 0.0833817738544132

$ws.Range("I9").Value = 0.2583188840501616
$ws.Range("J9").Value = 0.2583188840501616
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 16489.95503167943
$ws.Range("R9").Value = 148409.5952851149
$ws.Range("S9").Value = 0.03594478950200088
$ws.Range("T9").Value = 0.03594478950200088

$ws.Range("I10").Value = 0.2583188840501616
$ws.Range("J10").Value = 0.2583188840501616
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 12949.10704553697
$ws.Range("R10").Value = 116541.9634098327
$ws.Range("S10").Value = 0.02822645217021543
$ws.Range("T10").Value = 0.02822645217021544

$ws.Range("I11").Value = 0.2583188840501616
$ws.Range("J11").Value = 0.2583188840501616
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 15461.04904247865
$ws.Range("R11").Value = 139149.4413823079
$ws.Range("S11").Value = 0.0337019811299878
$ws.Range("T11").Value = 0.0337019811299878

$ws.Range("G12").Value = 242.2070976666666
$ws.Range("H12").Value = 726.6212929999999
$ws.Range("I12").Value = 0.2822439972492947
$ws.Range("J12").Value = 0.2822439972492947
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 38628.06892382991
$ws.Range("R12").Value = 347652.6203144692
$ws.Range("S12").Value = 0.08420143072970149
$ws.Range("T12").Value = 0.08420143072970147

$ws.Range("G13").Value = 242.2070976666666
$ws.Range("H13").Value = 726.6212929999999
$ws.Range("I13").Value = 0.2822439972492947
$ws.Range("J13").Value = 0.2822439972492947
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 41794.8927360922
$ws.Range("R13").Value = 376154.0346248298
$ws.Range("S13").Value = 0.09110447049560791
$ws.Range("T13").Value = 0.09110447049560789

$ws.Range("G14").Value = 242.2070976666666
$ws.Range("H14").Value = 726.6212929999999
$ws.Range("I14").Value = 0.2822439972492947
$ws.Range("J14").Value = 0.2822439972492947
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 18017.23029160559
$ws.Range("R14").Value = 162155.0726244503
$ws.Range("S14").Value = 0.0392739427728372
$ws.Range("T14").Value = 0.0392739427728372

$ws.Range("G15").Value = 242.2070976666666
$ws.Range("H15").Value = 726.6212929999999
$ws.Range("I15").Value = 0.2822439972492947
$ws.Range("J15").Value = 0.2822439972492947
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 14148.43419899433
$ws.Range("R15").Value = 127335.907790949
$ws.Range("S15").Value = 0.03084074444646723
$ws.Range("T15").Value = 0.03084074444646723

$ws.Range("G16").Value = 242.2070976666666
$ws.Range("H16").Value = 726.6212929999999
$ws.Range("I16").Value = 0.2822439972492947
$ws.Range("J16").Value = 0.2822439972492947
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 16893.02855059243
$ws.Range("R16").Value = 152037.2569553319
$ws.Range("S16").Value = 0.03682340880468087
$ws.Range("T16").Value = 0.03682340880468087

$ws.Range("G17").Value = 9.570116333333333
$ws.Range("H17").Value = 28.710349
$ws.Range("I17").Value = 0.01115205918440143
$ws.Range("J17").Value = 0.01115205918440143
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 1526.276962543143
$ws.Range("R17").Value = 13736.49266288829
$ws.Range("S17").Value = 0.003326977183077202
$ws.Range("T17").Value = 0.003326977183077202

$ws.Range("G18").Value = 9.570116333333333
$ws.Range("H18").Value = 28.710349
$ws.Range("I18").Value = 0.01115205918440143
$ws.Range("J18").Value = 0.01115205918440143
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 1651.404890595151
$ws.Range("R18").Value = 14862.64401535636
$ws.Range("S18").Value = 0.003599730930798785
$ws.Range("T18").Value = 0.003599730930798785

$ws.Range("G19").Value = 9.570116333333333
$ws.Range("H19").Value = 28.710349
$ws.Range("I19").Value = 0.01115205918440143
$ws.Range("J19").Value = 0.01115205918440143
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 711.8989969997591
$ws.Range("R19").Value = 6407.090972997833
$ws.Range("S19").Value = 0.001551796808704563
$ws.Range("T19").Value = 0.001551796808704564

$ws.Range("G20").Value = 9.570116333333333
$ws.Range("H20").Value = 28.710349
$ws.Range("I20").Value = 0.01115205918440143
$ws.Range("J20").Value = 0.01115205918440143
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 559.0346547368008
$ws.Range("R20").Value = 5031.311892631207
$ws.Range("S20").Value = 0.001218583249634946
$ws.Range("T20").Value = 0.001218583249634946

$ws.Range("G21").Value = 9.570116333333333
$ws.Range("H21").Value = 28.710349
$ws.Range("I21").Value = 0.01115205918440143
$ws.Range("J21").Value = 0.01115205918440143
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 667.4794009297947
$ws.Range("R21").Value = 6007.314608368152
$ws.Range("S21").Value = 0.001454971012185932
$ws.Range("T21").Value = 0.001454971012185932

$ws.Range("G22").Value = 315.9253133333333
$ws.Range("H22").Value = 947.77594
$ws.Range("I22").Value = 0.3681478541564123
$ws.Range("J22").Value = 0.3681478541564123
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 50384.91809607304
$ws.Range("R22").Value = 453464.2628646573
$ws.Range("S22").Value = 0.1098290002343597
$ws.Range("T22").Value = 0.1098290002343597

$ws.Range("G23").Value = 315.9253133333333
$ws.Range("H23").Value = 947.77594
$ws.Range("I23").Value = 0.3681478541564123
$ws.Range("J23").Value = 0.3681478541564123
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 54515.59723305406
$ws.Range("R23").Value = 490640.3750974865
$ws.Range("S23").Value = 0.1188330509909473
$ws.Range("T23").Value = 0.1188330509909473

$ws.Range("G24").Value = 315.9253133333333
$ws.Range("H24").Value = 947.77594
$ws.Range("I24").Value = 0.3681478541564123
$ws.Range("J24").Value = 0.3681478541564123
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 23500.95922088944
$ws.Range("R24").Value = 211508.632988005
$ws.Range("S24").Value = 0.05122737027888333
$ws.Range("T24").Value = 0.05122737027888334

$ws.Range("G25").Value = 315.9253133333333
$ws.Range("H25").Value = 947.77594
$ws.Range("I25").Value = 0.3681478541564123
$ws.Range("J25").Value = 0.3681478541564123
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 18454.65533650416
$ws.Range("R25").Value = 166091.8980285374
$ws.Range("S25").Value = 0.04022744150170433
$ws.Range("T25").Value = 0.04022744150170433

$ws.Range("G26").Value = 315.9253133333333
$ws.Range("H26").Value = 947.77594
$ws.Range("I26").Value = 0.3681478541564123
$ws.Range("J26").Value = 0.3681478541564123
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 22034.59514361435
$ws.Range("R26").Value = 198311.3562925291
$ws.Range("S26").Value = 0.04803099115051765
$ws.Range("T26").Value = 0.04803099115051765
